$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 140, pushing existing rows 140-174 down to 141-175.
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new record's data.
$ws.Cells.Item(140, 1).Value = 7
$ws.Cells.Item(140, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(140, 3).Value = "Ñuble"
$ws.Cells.Item(140, 4).Value = 44798
$ws.Cells.Item(140, 5).Value = 16
$ws.Cells.Item(140, 6).Value = 100112045
$ws.Cells.Item(140, 7).Value = "Zapallo"
$ws.Cells.Item(140, 8).Value = "Camote"
$ws.Cells.Item(140, 9).Value = "1a (guarda)"
$ws.Cells.Item(140, 10).Value = 300
$ws.Cells.Item(140, 11).Value = 800
$ws.Cells.Item(140, 12).Value = 900
$ws.Cells.Item(140, 13).Value = 850
$ws.Cells.Item(140, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(140, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(140, 16).Value = 850
$ws.Cells.Item(140, 17).Value = 1
$ws.Cells.Item(140, 18).Value = "Hortaliza"
